$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (2-10): a full 3x3 grid of Sending cluster x Target cluster
# combinations (ECs / FAPs / sCs), replacing the original 3 rows (2-4) which
# only covered Target cluster = sCs.
$data = @(
    @("ECs","Fn1","Itgb6","ECs",3,1,21.84976866666667,65.549306,0.05020018890879543,0.05020018890879543,1,0.3333333333333333,0.069456,0.208368,0.01627409789654661,0.01627409789654661,1.517597532512,13.658377792608,0.0008169627887268699,0.0008169627887268699),
    @("ECs","Fn1","Itgb6","FAPs",3,1,21.84976866666667,65.549306,0.05020018890879543,0.05020018890879543,2,0.6666666666666666,0.6957970000000001,2.087391,0.1630308179872645,0.1630308179872644,15.20300348896067,136.827031400646,0.008184177860916119,0.008184177860916119),
    @("ECs","Fn1","Itgb6","sCs",3,1,21.84976866666667,65.549306,0.05020018890879543,0.05020018890879543,3,1,3.502633333333333,10.5079,0.820695084116189,0.820695084116189,76.53172805748889,688.7855525174,0.04119904825915244,0.04119904825915244),
    @("FAPs","Fn1","Itgb6","ECs",3,1,385.0524703333334,1155.157411,0.8846641374295412,0.8846641374295412,1,0.3333333333333333,0.069456,0.208368,0.01627409789654661,0.01627409789654661,26.744204379472,240.697839415248,0.01439711077809231,0.01439711077809231),
    @("FAPs","Fn1","Itgb6","FAPs",3,1,385.0524703333334,1155.157411,0.8846641374295412,0.8846641374295412,2,0.6666666666666666,0.6957970000000001,2.087391,0.1630308179872645,0.1630308179872644,267.9183537005224,2411.265183304702,0.1442275179691359,0.1442275179691358),
    @("FAPs","Fn1","Itgb6","sCs",3,1,385.0524703333334,1155.157411,0.8846641374295412,0.8846641374295412,3,1,3.502633333333333,10.5079,0.820695084116189,0.820695084116189,1348.697617671878,12138.2785590469,0.726039508682313,0.726039508682313),
    @("sCs","Fn1","Itgb6","ECs",3,1,28.350479,85.05143699999999,0.06513567366166337,0.06513567366166337,1,0.3333333333333333,0.069456,0.208368,0.01627409789654661,0.01627409789654661,1.969110869424,17.721997824816,0.001060024329727422,0.001060024329727422),
    @("sCs","Fn1","Itgb6","FAPs",3,1,28.350479,85.05143699999999,0.06513567366166337,0.06513567366166337,2,0.6666666666666666,0.6957970000000001,2.087391,0.1630308179872645,0.1630308179872644,19.726178236763,177.535604130867,0.0106191221572125,0.01061912215721249),
    @("sCs","Fn1","Itgb6","sCs",3,1,28.350479,85.05143699999999,0.06513567366166337,0.06513567366166337,3,1,3.502633333333333,10.5079,0.820695084116189,0.820695084116189,99.30133276136665,893.7119948522999,0.05345652717472345,0.05345652717472345)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $c = $j + 1
        $ws.Cells.Item($r, $c).Value2 = $rowVals[$j]
    }
}
